$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '36.390.43'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.69%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.952.68'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -3.33%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '243.21'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -1.84%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.613'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -3.54%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '57.96'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -7.12%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.08%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.370'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -4.73%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '55.70'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -3.62%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0821'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +4.11%  '
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +0.10%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.830'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -7.69%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '21.56'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -7.05%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.237.96'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -3.72%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '13.53'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -5.33%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '5.29'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -4.09%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.945.70'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -3.82%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '36.281.21'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -0.81%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '69.87'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -3.09%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0₃0873'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -1.01%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '229.62'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -2.63%  '
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -6.50%  '
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +0.20%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.52'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +0.39%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.28'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -1.86%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.47'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -3.44%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '163.80'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +2.67%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '19.54'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -3.72%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.119'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -11.16%  '
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -1.83%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -1.43%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.73'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -6.30%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0635'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +2.78%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.32'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -3.50%  '
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +0.11%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.08'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -4.46%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -1.69%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.15'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -9.01%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.89'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -9.10%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0987'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -0.97%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.87'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -2.52%  '
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -4.69%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0210'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -2.10%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '15.70'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -7.54%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.04'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -8.11%  '
$ws.Range('B47').Value = 'FraxShare'
$ws.Range('C47').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '7.38'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -3.62%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.343.54'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -1.90%  '
$ws.Range('B49').Value = 'Aave'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '88.36'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -5.60%  '
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -2.34%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '45.85'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +1.67%  '
